$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column B (Part Code) / Column C (Part Name) / Column D (Part Number) mapping.
# The "Mechanical" block (rows 6-8) now reuses the same Electrical part-code
# series (ELP00111..117 / EL Part 111..117) instead of the old MCP/MC Part
# values, and column D switches from "Part noX" text codes to "numX" values.
$ws.Range("B2").Value = "ELP00111"
$ws.Range("C2").Value = "EL Part 111"
$ws.Range("D2").Value = "num1"

$ws.Range("B3").Value = "ELP00112"
$ws.Range("C3").Value = "EL Part 112"
$ws.Range("D3").Value = "num2"

$ws.Range("B4").Value = "ELP00113"
$ws.Range("C4").Value = "EL Part 113"
$ws.Range("D4").Value = "num3"

$ws.Range("B5").Value = "ELP00114"
$ws.Range("C5").Value = "EL Part 114"
$ws.Range("D5").Value = "num4"

$ws.Range("B6").Value = "ELP00115"
$ws.Range("C6").Value = "EL Part 115"
$ws.Range("D6").Value = "num5"

$ws.Range("B7").Value = "ELP00116"
$ws.Range("C7").Value = "EL Part 116"
$ws.Range("D7").Value = "num6"

$ws.Range("B8").Value = "ELP00117"
$ws.Range("C8").Value = "EL Part 117"
$ws.Range("D8").Value = "num76"

# Move the active selection to D8, matching the saved view state.
$ws.Range("D8").Select()
